$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (the "Identifier" column). This shifts all remaining
# columns (and their data) left by one, updates the sheet dimension, the
# column-width definitions, and drops the now-unused "Identifier" entry
# from the shared-string table.
$ws.Columns.Item(1).Delete()

# The conditional-formatting rules' applied ranges still reference the
# pre-delete column letters, so re-point them at the post-delete columns:
# old column C rules -> new column B, old column D rules -> new column C.
$newB = $ws.Range("B1:B1048576")
$newC = $ws.Range("C1:C1048576")

$cRange = $ws.Range("C1:C1048576")
$cCount = $cRange.FormatConditions.Count
for ($i = 1; $i -le $cCount; $i++) {
    $cRange.FormatConditions.Item(1).ModifyAppliesToRange($newB)
}

$dRange = $ws.Range("D1:D1048576")
$dCount = $dRange.FormatConditions.Count
for ($i = 1; $i -le $dCount; $i++) {
    $dRange.FormatConditions.Item(1).ModifyAppliesToRange($newC)
}

# Set the active selection to match the post-edit state.
$ws.Range("E10").Select()
